# Generate Report for Handoff
# Rename the localized file's GUID-based base name from
#   1274b2b4-be45-488d-a0c1-f5187d565386
# to
#   e414559a-85d2-4c60-8b29-5c9aa639a168
# across all three sheets, refresh the handoff/handback timestamps, and
# keep each hyperlinked cell's display text in sync with its new value.

$wb = $excel.ActiveWorkbook

$oldBase = "1274b2b4-be45-488d-a0c1-f5187d565386"
$newBase = "e414559a-85d2-4c60-8b29-5c9aa639a168"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9c56eeef66cc751ed73c220cf6381f1eb698af2/e2e/$oldBase.md"

function Set-HyperlinkCell($Sheet, $Cell, $NewText) {
    # Hyperlinks.Delete()/Add() operate on the sheet's whole hyperlink
    # collection in this host, so only call this for the cell that truly
    # owns the hyperlink on that sheet (one per sheet here).
    $range = $Sheet.Range($Cell)
    $range.Hyperlinks.Delete()
    $range.Value = $NewText
    $Sheet.Hyperlinks.Add($range, $hyperlinkTarget, "", "", $NewText) | Out-Null
}

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "$newBase.md"
Set-HyperlinkCell $ws "B2" "e2e\$newBase.md"
$ws.Range("G2").Value = "2016-09-06 21:18:33"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
Set-HyperlinkCell $ws "A2" "$newBase.md"
$ws.Range("G2").Value = "$newBase.3bfab2c72810c51a52d881e839aa9c9ddc79520a.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-06 21:18:27"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
Set-HyperlinkCell $ws "A2" "$newBase.md"
$ws.Range("G2").Value = "$newBase.3bfab2c72810c51a52d881e839aa9c9ddc79520a.de-de.xlf"
$ws.Range("H2").Value = "2016-09-06 21:18:33"
